$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 5.313941333333333
$ws.Range("H2").Value = 15.941824
$ws.Range("I2").Value = 0.176869630377001
$ws.Range("J2").Value = 0.176869630377001
$ws.Range("M2").Value = 3.456265333333333
$ws.Range("N2").Value = 10.368796
$ws.Range("O2").Value = 0.009841535807677501
$ws.Range("P2").Value = 0.0098415358076775
$ws.Range("Q2").Value = 18.36639121376711
$ws.Range("R2").Value = 165.297520923904
$ws.Range("S2").Value = 0.001740668800645939
$ws.Range("T2").Value = 0.001740668800645939
$ws.Range("G3").Value = 5.313941333333333
$ws.Range("H3").Value = 15.941824
$ws.Range("I3").Value = 0.176869630377001
$ws.Range("J3").Value = 0.176869630377001
$ws.Range("O3").Value = 0.8587907398420774
$ws.Range("P3").Value = 0.8587907398420773
$ws.Range("Q3").Value = 1602.685496139278
$ws.Range("R3").Value = 14424.1694652535
$ws.Range("S3").Value = 0.1518940007270594
$ws.Range("T3").Value = 0.1518940007270594
$ws.Range("G4").Value = 5.313941333333333
$ws.Range("H4").Value = 15.941824
$ws.Range("I4").Value = 0.176869630377001
$ws.Range("J4").Value = 0.176869630377001
$ws.Range("O4").Value = 0.1313677243502452
$ws.Range("P4").Value = 0.1313677243502452
$ws.Range("Q4").Value = 245.1600101273529
$ws.Range("R4").Value = 2206.440091146176
$ws.Range("S4").Value = 0.02323496084929562
$ws.Range("T4").Value = 0.02323496084929561
$ws.Range("I5").Value = 0.5461014638447835
$ws.Range("J5").Value = 0.5461014638447835
$ws.Range("M5").Value = 3.456265333333333
$ws.Range("N5").Value = 10.368796
$ws.Range("O5").Value = 0.009841535807677501
$ws.Range("P5").Value = 0.0098415358076775
$ws.Range("Q5").Value = 56.70794418468134
$ws.Range("R5").Value = 510.371497662132
$ws.Range("S5").Value = 0.005374477111053538
$ws.Range("T5").Value = 0.005374477111053537
$ws.Range("I6").Value = 0.5461014638447835
$ws.Range("J6").Value = 0.5461014638447835
$ws.Range("O6").Value = 0.8587907398420774
$ws.Range("P6").Value = 0.8587907398420773
$ws.Range("S6").Value = 0.4689868801641032
$ws.Range("T6").Value = 0.4689868801641031
$ws.Range("I7").Value = 0.5461014638447835
$ws.Range("J7").Value = 0.5461014638447835
$ws.Range("O7").Value = 0.1313677243502452
$ws.Range("P7").Value = 0.1313677243502452
$ws.Range("S7").Value = 0.07174010656962691
$ws.Range("T7").Value = 0.0717401065696269
$ws.Range("I8").Value = 0.2770289057782155
$ws.Range("J8").Value = 0.2770289057782155
$ws.Range("M8").Value = 3.456265333333333
$ws.Range("N8").Value = 10.368796
$ws.Range("O8").Value = 0.009841535807677501
$ws.Range("P8").Value = 0.0098415358076775
$ws.Range("Q8").Value = 28.76707126146711
$ws.Range("R8").Value = 258.903641353204
$ws.Range("S8").Value = 0.002726389895978025
$ws.Range("T8").Value = 0.002726389895978024
$ws.Range("I9").Value = 0.2770289057782155
$ws.Range("J9").Value = 0.2770289057782155
$ws.Range("O9").Value = 0.8587907398420774
$ws.Range("P9").Value = 0.8587907398420773
$ws.Range("S9").Value = 0.2379098589509149
$ws.Range("T9").Value = 0.2379098589509148
$ws.Range("I10").Value = 0.2770289057782155
$ws.Range("J10").Value = 0.2770289057782155
$ws.Range("O10").Value = 0.1313677243502452
$ws.Range("P10").Value = 0.1313677243502452
$ws.Range("S10").Value = 0.03639265693132267
$ws.Range("T10").Value = 0.03639265693132266
